$d = $word.ActiveDocument

$replacements = @(
    @("2025-12-30 Tuesday", "2025-12-31 Wednesday"),
    @("97÷2=", "83÷3="),
    @("53÷2=", "14÷9="),
    @("16÷8=", "48÷3="),
    @("71÷7=", "38÷8="),
    @("26÷3=", "39÷5="),
    @("84÷9=", "78÷6="),
    @("31÷6=", "64÷7="),
    @("90÷9=", "65÷7="),
    @("15÷9=", "35÷4="),
    @("14÷7=", "93÷5="),
    @("55÷3=", "28÷7="),
    @("85÷4=", "60÷8="),
    @("28÷4=", "35÷9="),
    @("29÷2=", "90÷9="),
    @("60÷2=", "77÷8="),
    @("88÷7=", "88÷8="),
    @("24÷6=", "32÷4="),
    @("39÷2=", "65÷5="),
    @("27÷3=", "96÷7="),
    @("14÷8=", "25÷4="),
    @("93÷6=", "33÷9="),
    @("39÷3=", "30÷6="),
    @("91÷6=", "83÷4="),
    @("76÷8=", "19÷6="),
    @("20÷3=", "86÷6=")
)

foreach ($pair in $replacements) {
    $old = $pair[0]
    $new = $pair[1]
    $range = $d.Content
    $range.Find.Execute($old, $true, $true, $false, $false, $false, $true, 1, $false, $new, 2)
}
